# Updated results for BoW+TFIDF Model - Lucene, Thunderbird and Ubuntu
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Lucene sheet: add new "Count Vectorizer + TFIDF" block (rows 9-13)
# ---------------------------------------------------------------
$wsLucene = $wb.Worksheets.Item("Lucene")

$wsLucene.Range("A9").Value = "Logistic Regression"
$wsLucene.Range("B9").Value = "Count Vectorizer + TFIDF"
$wsLucene.Range("C9").Value = "0.377 0.657 0.512 0.368 0.743"
$wsLucene.Range("D9").Value = "0.123 0.674 0.686 0.132 0.879"
$wsLucene.Range("E9").Value = "0.247 0.528 0.350 0.239 0.594"
$wsLucene.Range("F9").Value = "0.733 0.754 0.820 0.729 0.959"

$wsLucene.Range("A10").Value = "Multinomial Naive Bayes"
$wsLucene.Range("B10").Value = "Count Vectorizer + TFIDF"
$wsLucene.Range("C10").Value = "0.416 0.740 0.719 0.609 0.845"
$wsLucene.Range("D10").Value = "0.383 0.584 0.501 0.321 0.606"
$wsLucene.Range("E10").Value = "0.266 0.755 0.641 0.468 0.758"
$wsLucene.Range("F10").Value = "0.877 0.735 0.778 0.823 0.938"

$wsLucene.Range("A11").Value = "Support Vector Machines"
$wsLucene.Range("B11").Value = "Count Vectorizer + TFIDF"
$wsLucene.Range("C11").Value = "0.417 0.725 0.694 0.522 0.786"
$wsLucene.Range("D11").Value = "0.437 0.699 0.671 0.507 0.839"
$wsLucene.Range("E11").Value = "0.266 0.626 0.556 0.359 0.652"
$wsLucene.Range("F11").Value = "0.886 0.782 0.841 0.885 0.960"

$wsLucene.Range("A12").Value = "Decision Tree"
$wsLucene.Range("B12").Value = "Count Vectorizer + TFIDF"
$wsLucene.Range("C12").Value = "0.060 0.608 0.495 0.162 0.667"
$wsLucene.Range("D12").Value = "0.571 0.618 0.693 0.391 0.852"
$wsLucene.Range("E12").Value = "0.031 0.473 0.333 0.088 0.502"
$wsLucene.Range("F12").Value = "0.895 0.723 0.819 0.878 0.951"

$wsLucene.Range("A13").Value = "Random Forest"
$wsLucene.Range("B13").Value = "Count Vectorizer + TFIDF"
$wsLucene.Range("C13").Value = "0.176 0.673 0.461 0.252 0.681"
$wsLucene.Range("D13").Value = "1.000 0.831 0.911 0.911 0.947"
$wsLucene.Range("E13").Value = "0.097 0.523 0.300 0.144 0.517"
$wsLucene.Range("F13").Value = "0.904 0.803 0.838 0.899 0.957"

# ---------------------------------------------------------------
# Thunderbird sheet: add new "Count Vectorizer + TFIDF" block (rows 9-13)
# ---------------------------------------------------------------
$wsThunderbird = $wb.Worksheets.Item("Thunderbird")

$wsThunderbird.Range("A9").Value = "Logistic Regression"
$wsThunderbird.Range("B9").Value = "Count Vectorizer + TFIDF"
$wsThunderbird.Range("C9").Value = "0.477 0.397 0.338 0.494 0.416"
$wsThunderbird.Range("D9").Value = "0.613 0.602 0.076 0.077 0.677"
$wsThunderbird.Range("E9").Value = "0.320 0.250 0.214 0.421 0.263"
$wsThunderbird.Range("F9").Value = "0.786 0.856 0.755 0.584 0.978"

$wsThunderbird.Range("A10").Value = "Multinomial Naive Bayes"
$wsThunderbird.Range("B10").Value = "Count Vectorizer + TFIDF"
$wsThunderbird.Range("C10").Value = "0.669 0.681 0.512 0.550 0.592"
$wsThunderbird.Range("D10").Value = "0.491 0.479 0.292 0.309 0.293"
$wsThunderbird.Range("E10").Value = "0.569 0.552 0.353 0.390 0.425"
$wsThunderbird.Range("F10").Value = "0.753 0.836 0.891 0.890 0.959"

$wsThunderbird.Range("A11").Value = "Support Vector Machines"
$wsThunderbird.Range("B11").Value = "Count Vectorizer + TFIDF"
$wsThunderbird.Range("C11").Value = "0.642 0.613 0.415 0.483 0.490"
$wsThunderbird.Range("D11").Value = "0.614 0.587 0.476 0.566 0.788"
$wsThunderbird.Range("E11").Value = "0.499 0.455 0.263 0.320 0.325"
$wsThunderbird.Range("F11").Value = "0.803 0.864 0.925 0.932 0.980"

$wsThunderbird.Range("A12").Value = "Decision Tree"
$wsThunderbird.Range("B12").Value = "Count Vectorizer + TFIDF"
$wsThunderbird.Range("C12").Value = "0.368 0.236 0.061 0.026 0.242"
$wsThunderbird.Range("D12").Value = "0.526 0.436 0.438 0.188 0.917"
$wsThunderbird.Range("E12").Value = "0.229 0.134 0.031 0.013 0.138"
$wsThunderbird.Range("F12").Value = "0.763 0.837 0.927 0.923 0.977"

$wsThunderbird.Range("A13").Value = "Random Forest"
$wsThunderbird.Range("B13").Value = "Count Vectorizer + TFIDF"
$wsThunderbird.Range("C13").Value = "0.315 0.302 0.229 0.266 0.462"
$wsThunderbird.Range("D13").Value = "0.892 0.935 0.879 0.921 1.000"
$wsThunderbird.Range("E13").Value = "0.187 0.178 0.129 0.154 0.300"
$wsThunderbird.Range("F13").Value = "0.798 0.869 0.936 0.936 0.982"

# ---------------------------------------------------------------
# Ubuntu sheet: values for its "Count Vectorizer + TFIDF" block (rows 9-13)
# already existed and are unchanged in content.
# ---------------------------------------------------------------
$wsUbuntu = $wb.Worksheets.Item("Ubuntu")

# ---------------------------------------------------------------
# Update sheet tab/view selection state:
#  - Lucene becomes the active (selected) tab, cell E13 selected
#  - Thunderbird and Ubuntu become inactive tabs, cell F13 selected
# ---------------------------------------------------------------
$wsThunderbird.Activate()
$wsThunderbird.Range("F13").Select()

$wsUbuntu.Activate()
$wsUbuntu.Range("F13").Select()

$wsLucene.Activate()
$wsLucene.Range("E13").Select()
